$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column E
$ws.Range("E1").Value = "Debug"

# Row 2: update Set point and Time (Debug left blank)
$ws.Range("C2").Value = 76.09999999999999
$ws.Range("D2").Value = "2016-06-02 15:26:11"

# Row 3
$ws.Range("A3").Value = 256.56
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 73.90000000000001
$ws.Range("D3").Value = "2016-06-03 10:02:14"
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("A4").Value = 256.56
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 73.90000000000001
$ws.Range("D4").Value = "2016-06-03 10:03:40"
$ws.Range("E4").Value = 1

# Row 5
$ws.Range("A5").Value = 256.56
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 73.90000000000001
$ws.Range("D5").Value = "2016-06-03 10:20:10"
$ws.Range("E5").Value = 1

# Row 6
$ws.Range("A6").Value = 256.56
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 73.90000000000001
$ws.Range("D6").Value = "2016-06-03 10:21:15"
$ws.Range("E6").Value = 0

# Row 7
$ws.Range("A7").Value = 71.03999999999999
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 71.7
$ws.Range("D7").Value = "2016-06-03 21:08:22"
$ws.Range("E7").Value = 1

# Row 8
$ws.Range("A8").Value = 71.03999999999999
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 71.7
$ws.Range("D8").Value = "2016-06-03 21:44:34"
$ws.Range("E8").Value = 1

# Row 9
$ws.Range("A9").Value = 71.03999999999999
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 71.7
$ws.Range("D9").Value = "2016-06-03 21:44:44"
$ws.Range("E9").Value = 1

# Row 10
$ws.Range("A10").Value = 71.03999999999999
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 71.7
$ws.Range("D10").Value = "2016-06-03 21:44:54"
$ws.Range("E10").Value = 1

# Row 11
$ws.Range("A11").Value = 71.03999999999999
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 71.7
$ws.Range("D11").Value = "2016-06-03 22:33:49"
$ws.Range("E11").Value = 1

# Row 12
$ws.Range("A12").Value = 71.03999999999999
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 71.7
$ws.Range("D12").Value = "2016-06-03 22:33:59"
$ws.Range("E12").Value = 1

# Row 13
$ws.Range("A13").Value = 71.03999999999999
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 71.7
$ws.Range("D13").Value = "2016-06-03 22:34:09"
$ws.Range("E13").Value = 1

# Row 14
$ws.Range("A14").Value = 71.03999999999999
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 71.7
$ws.Range("D14").Value = "2016-06-03 22:34:19"
$ws.Range("E14").Value = 1

# Row 15
$ws.Range("A15").Value = 71.03999999999999
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 71.7
$ws.Range("D15").Value = "2016-06-03 22:35:36"
$ws.Range("E15").Value = 0
